$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Index" -> "i" (column A). This also renames the table's
# ListColumn (the worksheet is backed by table "testdata").
$ws.Cells.Item(1, 1).Value2 = "i"

# Narrow column A (was width 6, now width 4).
$ws.Columns.Item(1).ColumnWidth = 3.14

# Re-index data rows from 1-based to 0-based (row 2 -> 0, row 3 -> 1, ...).
for ($r = 2; $r -le 503; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}
